# Auto-generated edit script: updates cryptocurrency price/volume/name/link
# cells per the "Updated cryptos list on Sat Oct 26 02:03:58 UTC 2024 with
# GitHub Actions" commit. Column D/E store numbers and percentages as plain
# text (European-style "66.666.93" style prices, "  -1.95%  " deltas), so
# any replacement value that Excel would otherwise auto-parse into a real
# number is written with a leading apostrophe (forces text) and the cell
# style is put back to Normal right after so no stray "quote prefix" / text
# number-format is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.666.93'
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").Value = '2.440.62'
$ws.Range("E3").Value = '  -3.16%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'577.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.94%  '
$ws.Range("D6").Value = "'163.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.24%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = "'0.508"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.24%  '
$ws.Range("D9").Value = '2.442.31'
$ws.Range("E9").Value = '  -2.96%  '
$ws.Range("D10").Value = "'0.132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.54%  '
$ws.Range("E11").Value = '  -1.15%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = "'0.329"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.35%  '
$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").Value = "'4.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.53%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.888.68'
$ws.Range("E14").Value = '  -3.00%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = "'24.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.04%  '
$ws.Range("D16").Value = '66.615.08'
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("E17").Value = '  -6.59%  '
$ws.Range("D18").Value = '2.440.69'
$ws.Range("E18").Value = '  -2.75%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = "'7.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.36%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'351.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.98%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = "'10.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -13.98%  '
$ws.Range("D22").Value = "'4.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.10%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = "'69.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'4.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -10.07%  '
$ws.Range("D26").Value = "'1.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -10.88%  '
$ws.Range("D27").Value = "'8.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -12.31%  '
$ws.Range("D28").Value = "'0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").Value = '2.564.75'
$ws.Range("E29").Value = '  -2.98%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0884'
$ws.Range("E30").Value = '  -9.91%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = "'502.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.65%  '
$ws.Range("D32").Value = "'7.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.73%  '
$ws.Range("D33").Value = "'1.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.09%  '
$ws.Range("E34").Value = '  -9.28%  '
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").Value = "'157.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("D37").Value = "'0.115"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -10.40%  '
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("D39").Value = "'18.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.77%  '
$ws.Range("E40").Value = '  -8.42%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'1.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.40%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").Value = "'0.323"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.31%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D44").Value = "'4.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.97%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = "'38.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.37%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = "'2.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.97%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = "'140.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.89%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").Value = "'3.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.05%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = "'0.509"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.37%  '
$ws.Range("B50").Value = 'Optimism'
$ws.Range("C50").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D50").Value = "'1.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.53%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.0726"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.68%  '
